# DICHIARATI I TIPI DI DATI
# Fill in the "Tipo di Dato :" row (row 4) with the SQL data type declared
# for each column described in row 3 of the "Database First" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "VARCHAR(100)"   # modello
$ws.Range("C4").Value = "TINYINT"        # porte
$ws.Range("D4").Value = "TINYINT"        # numero_Sedie
$ws.Range("E4").Value = "CHAR(7)"        # targa
$ws.Range("F4").Value = "TINYINT"        # numero_proprietari_precendti
$ws.Range("G4").Value = "SMALLINT"       # cilindrata
$ws.Range("H4").Value = "BOOL"           # uso_commerciale
$ws.Range("I4").Value = "BOOL"           # motore_originale

# Column B now needs to be widened to fit the newly typed values.
$ws.Columns("B:B").ColumnWidth = 12.8

# Move the active selection, as it was left after the edits.
$ws.Range("F7").Select() | Out-Null
